# Auto-generated Excel COM-interop script to apply profit-table updates
# across the Chocobo Profits workbook (recalculated leve-crafting values).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1969.1428
$ws.Range("I40").Value = 2326.6667
$ws.Range("J40").Value = 1701
$ws.Range("K40").Value = 2326.6667
$ws.Range("L40").Value = 1701
$ws.Range("M40").Value = -2151.6667
$ws.Range("N40").Value = -2051
$ws.Range("H70").Value = 2592
$ws.Range("I70").Value = 1471.25
$ws.Range("K70").Value = 4413.75
$ws.Range("M70").Value = -4143.75
$ws.Range("H73").Value = 2592
$ws.Range("I73").Value = 1471.25
$ws.Range("K73").Value = 4413.75
$ws.Range("M73").Value = -3477.75
$ws.Range("H74").Value = 6399.294
$ws.Range("I74").Value = 5264.5557
$ws.Range("K74").Value = 5264.5557
$ws.Range("M74").Value = -4328.5557
$ws.Range("H77").Value = 6399.294
$ws.Range("I77").Value = 5264.5557
$ws.Range("K77").Value = 26322.7785
$ws.Range("M77").Value = -21642.7785
$ws.Range("H126").Value = 41852.223
$ws.Range("J126").Value = 41852.223
$ws.Range("L126").Value = 41852.223
$ws.Range("N126").Value = -51732.223
$ws.Range("H135").Value = 748.4074000000001
$ws.Range("I135").Value = 635.58826
$ws.Range("J135").Value = 940.2
$ws.Range("K135").Value = 5720.29434
$ws.Range("L135").Value = 8461.800000000001
$ws.Range("M135").Value = -3185.29434
$ws.Range("N135").Value = -13531.8
$ws.Range("H137").Value = 629396.3
$ws.Range("I137").Value = 1987791.2
$ws.Range("J137").Value = 2444.8076
$ws.Range("K137").Value = 5963373.6
$ws.Range("L137").Value = 7334.4228
$ws.Range("M137").Value = -5960823.6
$ws.Range("N137").Value = -12434.4228
$ws.Range("H138").Value = 2539.1875
$ws.Range("I138").Value = 1559.1111
$ws.Range("J138").Value = 3799.2856
$ws.Range("K138").Value = 4677.3333
$ws.Range("L138").Value = 11397.8568
$ws.Range("M138").Value = 462.6666999999998
$ws.Range("N138").Value = -21677.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2146.125
$ws.Range("I61").Value = 2308.9092
$ws.Range("K61").Value = 2308.9092
$ws.Range("M61").Value = -2096.9092
$ws.Range("H112").Value = 34850
$ws.Range("J112").Value = 34850
$ws.Range("L112").Value = 34850
$ws.Range("N112").Value = -37804
$ws.Range("H118").Value = 28390
$ws.Range("J118").Value = 28390
$ws.Range("L118").Value = 28390
$ws.Range("N118").Value = -31704
$ws.Range("H132").Value = 3122.1936
$ws.Range("I132").Value = 2377.2173
$ws.Range("K132").Value = 7131.651899999999
$ws.Range("M132").Value = -4601.651899999999
$ws.Range("H136").Value = 2146.125
$ws.Range("I136").Value = 2308.9092
$ws.Range("K136").Value = 6926.7276
$ws.Range("M136").Value = -4376.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2421
$ws.Range("I86").Value = 2440
$ws.Range("J86").Value = 2345
$ws.Range("K86").Value = 2440
$ws.Range("L86").Value = 2345
$ws.Range("M86").Value = -1317
$ws.Range("N86").Value = -4591
$ws.Range("H89").Value = 2421
$ws.Range("I89").Value = 2440
$ws.Range("J89").Value = 2345
$ws.Range("K89").Value = 12200
$ws.Range("L89").Value = 11725
$ws.Range("M89").Value = -6584
$ws.Range("N89").Value = -22957
$ws.Range("H134").Value = 4309.3335
$ws.Range("I134").Value = 1572.2
$ws.Range("J134").Value = 5499.391
$ws.Range("K134").Value = 4716.6
$ws.Range("L134").Value = 16498.173
$ws.Range("M134").Value = -2181.6
$ws.Range("N134").Value = -21568.173

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5557452
$ws.Range("I16").Value = 11112560
$ws.Range("J16").Value = 2343.9
$ws.Range("K16").Value = 11112560
$ws.Range("L16").Value = 2343.9
$ws.Range("M16").Value = -11112273
$ws.Range("N16").Value = -2917.9
$ws.Range("H31").Value = 296736.66
$ws.Range("I31").Value = 1690764
$ws.Range("J31").Value = 3257.2104
$ws.Range("K31").Value = 1690764
$ws.Range("L31").Value = 3257.2104
$ws.Range("M31").Value = -1690469
$ws.Range("N31").Value = -3847.2104
$ws.Range("H34").Value = 296736.66
$ws.Range("I34").Value = 1690764
$ws.Range("J34").Value = 3257.2104
$ws.Range("K34").Value = 1690764
$ws.Range("L34").Value = 3257.2104
$ws.Range("M34").Value = -1690562
$ws.Range("N34").Value = -3661.2104
$ws.Range("H58").Value = 2553.9429
$ws.Range("I58").Value = 1459.75
$ws.Range("J58").Value = 4941.273
$ws.Range("K58").Value = 1459.75
$ws.Range("L58").Value = 4941.273
$ws.Range("M58").Value = -1256.75
$ws.Range("N58").Value = -5347.273
$ws.Range("H99").Value = 3981.8333
$ws.Range("I99").Value = 1788.2222
$ws.Range("K99").Value = 1788.2222
$ws.Range("M99").Value = -290.2221999999999
$ws.Range("H113").Value = 5557452
$ws.Range("I113").Value = 11112560
$ws.Range("J113").Value = 2343.9
$ws.Range("K113").Value = 11112560
$ws.Range("L113").Value = 2343.9
$ws.Range("M113").Value = -11110390
$ws.Range("N113").Value = -6683.9
$ws.Range("H126").Value = 3981.8333
$ws.Range("I126").Value = 1788.2222
$ws.Range("K126").Value = 5364.6666
$ws.Range("M126").Value = -2894.6666
$ws.Range("H136").Value = 2553.9429
$ws.Range("I136").Value = 1459.75
$ws.Range("J136").Value = 4941.273
$ws.Range("K136").Value = 4379.25
$ws.Range("L136").Value = 14823.819
$ws.Range("M136").Value = -1829.25
$ws.Range("N136").Value = -19923.819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2232212.5
$ws.Range("I2").Value = 51.6
$ws.Range("J2").Value = 3246831.2
$ws.Range("K2").Value = 309.6
$ws.Range("L2").Value = 19480987.2
$ws.Range("M2").Value = -196.6
$ws.Range("N2").Value = -19481213.2
$ws.Range("H24").Value = 1730
$ws.Range("I24").Value = 912.5
$ws.Range("J24").Value = 5000
$ws.Range("K24").Value = 2737.5
$ws.Range("L24").Value = 15000
$ws.Range("N24").Value = -15460
$ws.Range("M24").Value = -2507.5
$ws.Range("H34").Value = 10541.218
$ws.Range("J34").Value = 8225.200000000001
$ws.Range("L34").Value = 24675.6
$ws.Range("N34").Value = -24843.6
$ws.Range("H68").Value = 3151.0908
$ws.Range("I68").Value = 1063
$ws.Range("K68").Value = 3189
$ws.Range("M68").Value = -2378
$ws.Range("H71").Value = 3151.0908
$ws.Range("I71").Value = 1063
$ws.Range("K71").Value = 9567
$ws.Range("M71").Value = -5511
$ws.Range("H121").Value = 1950.8
$ws.Range("I121").Value = 914.8333
$ws.Range("J121").Value = 2065.9075
$ws.Range("K121").Value = 2744.4999
$ws.Range("L121").Value = 6197.7225
$ws.Range("M121").Value = -1434.4999
$ws.Range("N121").Value = -8817.7225

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2696.25
$ws.Range("I132").Value = 2076.054
$ws.Range("J132").Value = 4782.364
$ws.Range("K132").Value = 6228.162
$ws.Range("L132").Value = 14347.092
$ws.Range("M132").Value = -3698.162
$ws.Range("N132").Value = -19407.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1258.2727
$ws.Range("I16").Value = 1258.2727
$ws.Range("K16").Value = 1258.2727
$ws.Range("M16").Value = -1088.2727
$ws.Range("H132").Value = 3413.1052
$ws.Range("I132").Value = 2461.9312
$ws.Range("J132").Value = 6478
$ws.Range("K132").Value = 7385.7936
$ws.Range("L132").Value = 19434
$ws.Range("M132").Value = -4855.7936
$ws.Range("N132").Value = -24494
$ws.Range("H136").Value = 4886.48
$ws.Range("I136").Value = 2877.9092
$ws.Range("K136").Value = 8633.7276
$ws.Range("M136").Value = -6083.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2242.1396
$ws.Range("I132").Value = 1282.4642
$ws.Range("J132").Value = 4033.5334
$ws.Range("K132").Value = 3847.3926
$ws.Range("L132").Value = 12100.6002
$ws.Range("M132").Value = -1317.3926
$ws.Range("N132").Value = -17160.6002
$ws.Range("H136").Value = 2615.3635
$ws.Range("I136").Value = 1088.5927
$ws.Range("J136").Value = 5040.2354
$ws.Range("K136").Value = 3265.7781
$ws.Range("L136").Value = 15120.7062
$ws.Range("M136").Value = -715.7780999999995
$ws.Range("N136").Value = -20220.7062
